$d = $word.ActiveDocument

# 1. Update the ID marker paragraph's text, removing the trailing space run
#    by including it in the search string and leaving it out of the
#    replacement (the merged run keeps the first run's formatting).
$d.Content.Find.Execute(
    "**ID__AFFARS_5316_topic_18__ID** ", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_5316_504__ID**", 2
)

# 2. Bring that paragraph's formatting in line with the rest of the body
#    text paragraphs: a zero-width paragraph border (used purely to carry
#    the w:space="5" spacing on all four sides) and a 225-twip left indent.
$p = $d.Paragraphs(1)
$p.Format.LeftIndent = 11.25

$borders = $p.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
